# Saldo.xlsx update
#
# Changes applied (see commit diff):
#   1. Row with Conta 005662526 (AGUINALDO): Saldo 80873.67 -> 80757.05
#   2. Row with Conta 004491730 (DENISE) removed entirely
#   3. New row inserted: Conta 001761119 / Nome BLUEMETRIX / Saldo 259.41
#      (placed immediately before the existing 004363260 / LARISSA row)
#   4. New row inserted: Conta 004273239 / Nome DANIEL / Saldo -142.79
#      (placed immediately before the existing 005135281 / RAFAEL row,
#       i.e. right after the 004976625 / NORTON row)
#
# Operations are performed from the bottom of the sheet upward so that
# earlier (lower) row numbers used below stay valid after each insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (preserving any leading zeros)
# instead of letting Excel auto-coerce a numeric-looking string to a
# number, then drop the temporary "Text" number format so the cell is
# left without any extra style applied (matching the rest of the sheet).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Locate the anchor rows by scanning column A so the script is resilient
# even if row numbers were ever slightly different than expected.
function Find-RowByAccount($account) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value() -eq $account) {
            return $r
        }
    }
    return -1
}

# --- 4. Insert DANIEL row right before RAFAEL (after NORTON) ---
$rowRafael = Find-RowByAccount "005135281"
$ws.Rows.Item($rowRafael).Insert()
Set-TextValue $ws.Cells.Item($rowRafael, 1) "004273239"
Set-TextValue $ws.Cells.Item($rowRafael, 2) "DANIEL"
$ws.Cells.Item($rowRafael, 3).Value = -142.79

# --- 3. Insert BLUEMETRIX row right before LARISSA (004363260) ---
$rowLarissa = Find-RowByAccount "004363260"
$ws.Rows.Item($rowLarissa).Insert()
Set-TextValue $ws.Cells.Item($rowLarissa, 1) "001761119"
Set-TextValue $ws.Cells.Item($rowLarissa, 2) "BLUEMETRIX"
$ws.Cells.Item($rowLarissa, 3).Value = 259.41

# --- 2. Delete the DENISE row entirely ---
$rowDenise = Find-RowByAccount "004491730"
if ($rowDenise -ne -1) {
    $ws.Rows.Item($rowDenise).Delete()
}

# --- 1. Update AGUINALDO's balance ---
$rowAguinaldo = Find-RowByAccount "005662526"
$ws.Cells.Item($rowAguinaldo, 3).Value = 80757.05

Write-Host "Saldo.xlsx updated."
